$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (MAE)
$ws.Range("B2").Value = 1.686
$ws.Range("C2").Value = 1.425
$ws.Range("D2").Value = 1.62
$ws.Range("E2").Value = 0.989
$ws.Range("F2").Value = 2.296

# Row 3 (MSE)
$ws.Range("B3").Value = 3.944
$ws.Range("C3").Value = 3.817
$ws.Range("D3").Value = 4.436
$ws.Range("E3").Value = 1.909
$ws.Range("F3").Value = 7.608

# Row 4 (mean Y-Test)
$ws.Range("B4").Value = 18.203
$ws.Range("C4").Value = 15.45
$ws.Range("D4").Value = 18.059
$ws.Range("E4").Value = 12.974
$ws.Range("F4").Value = 30.217

# Row 5 (mean Y-predicted)
$ws.Range("B5").Value = 18.253
$ws.Range("C5").Value = 15.18
$ws.Range("D5").Value = 17.634
$ws.Range("E5").Value = 12.964
$ws.Range("F5").Value = 30.592

# Row 6 (R2)
$ws.Range("B6").Value = 0.671
$ws.Range("C6").Value = 0.761
$ws.Range("D6").Value = 0.811
$ws.Range("E6").Value = 0.431
$ws.Range("F6").Value = 0.806
